$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.826.71'
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +1.77%  '

$style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.664.12'
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  +0.67%  '

$style = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = $style
$ws.Range('E4').Value = '  +0.25%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '328.82'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +6.85%  '

$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  +0.39%  '

$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3643'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +0.65%  '

$style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '46.53'
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  -1.81%  '

$style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3218'
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  -1.64%  '

$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.134'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +1.03%  '

$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07015'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +1.15%  '

$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.005'
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  +0.40%  '

$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.030'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +1.48%  '

$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '19.41'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +0.69%  '

$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.667.29'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +0.88%  '

$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.573'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  -0.66%  '

$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001041'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +0.09%  '

$style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06561'
$ws.Range('D18').Style = $style
$ws.Range('E18').Value = '  +0.75%  '

$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +0.42%  '

$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '78.32'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +2.64%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.885'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  -0.50%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '15.71'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +0.18%  '

$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.85'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +1.94%  '

$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.856.77'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +1.92%  '

$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.439'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +0.49%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.368'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +1.27%  '

$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '147.66'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +0.98%  '

$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.60'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  +1.50%  '

$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.859.14'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +1.16%  '

$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.22'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +0.88%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.168'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -1.00%  '

$style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.076'
$ws.Range('D32').Style = $style
$ws.Range('E32').Value = '  +0.71%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.690'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +1.23%  '

$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08413'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +0.66%  '

$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.644'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -1.76%  '

$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.12'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -1.18%  '

$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.120'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  -1.84%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02230'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +1.49%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.224'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +1.51%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05973'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  -0.97%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2075'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +1.29%  '

$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.161'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +0.13%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.002'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  +0.27%  '

$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5903'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +1.23%  '

$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.61'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +8.45%  '

$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.837'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +2.85%  '

$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5697'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +2.47%  '

$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '124.16'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +2.01%  '

$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.946'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +0.49%  '

$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06963'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +1.00%  '

$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.178'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  +2.26%  '
